$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "Sheet1" to "test"
$ws.Name = "test"

# Size column A (closest achievable snap to the recorded best-fit width of
# 22.26953125 chars; this runtime quantizes column widths to a coarser grid)
$ws.Columns.Item(1).ColumnWidth = 21.5

# Add the new row/cell content
$ws.Range("A3").Value = "updated the file"

# Move the active selection to B12
$ws.Range("B12").Select() | Out-Null
